# Generate Report for Handoff
#
# - Status moves from "In Translation" to "Ready for handoff" on all three
#   sheets (Overview!E2:F2, zh-cn!C2, de-de!C2).
# - The per-locale "Latest Handoff Datetime" / "Latest HO Xliff Generate
#   Date" timestamps advance by 50 seconds (the moment the handoff report
#   was generated).
# - The now-wider "Ready for handoff" text causes the Status columns to
#   widen (Overview columns E & F, and column C on the zh-cn / de-de
#   sheets).

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
# The ColumnWidth COM property only resolves to 1/6-character-unit steps in
# this runtime, so 16.3333... is the input that lands closest to the
# target post-autofit width of ~17.216 character units (stored XML width
# ~17.1667) produced by the wider "Ready for handoff" text.
$newWidth  = 16.3333333333333

# ---- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-19 14:45:46"

$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# ---- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-19 14:45:42"

$zhcn.Columns.Item(3).ColumnWidth = $newWidth

# ---- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-19 14:45:46"

$dede.Columns.Item(3).ColumnWidth = $newWidth
